# Feb 19 - results shown as a list view
#
# Remove the three rows that were resolved/dropped from the features list
# ("Be more efficient with location listener...", "Re-factor code into
# code blocks", "Get rid of Toast that comes up when location is found"),
# which were rows 5-7, shifting the remaining rows up, then append a new
# row describing the follow-up task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old rows 5, 6 and 7 (entire rows), shifting rows 8-10 up to 5-7.
$ws.Range("A5:D7").EntireRow.Delete() | Out-Null

# Add the new row (now row 8) describing the follow-up task.
$ws.Range("A8").Value = "Location listener called when orientation changes"
$ws.Range("B8").Value = "Bug"
$ws.Range("C8").Value = "Medium"
$ws.Range("D8").Value = "Developer"

# Update the selected/active cell to A9, matching the saved view state.
$ws.Range("A9").Select() | Out-Null
